# Append 9 new rows (102-110) of vehicle rental log data to Sheet1,
# mirroring the structure of existing rows in the table (A:N).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 102
$ws.Range("A102").Value = "'01-03-2018"
$ws.Range("A102").Style = "Normal"
$ws.Range("B102").Value = 'v4'
$ws.Range("C102").Value = 'Day/'
$ws.Range("D102").Value = 2000
$ws.Range("E102").Value = 2016
$ws.Range("F102").Value = 1700
$ws.Range("G102").Value = 16
$ws.Range("H102").Value = 66
$ws.Range("I102").Value = 200
$ws.Range("J102").Value = 27200
$ws.Range("K102").Value = 'Nothing'
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 14000

# Row 103
$ws.Range("A103").Value = "'01-03-2018"
$ws.Range("A103").Style = "Normal"
$ws.Range("B103").Value = 'v4'
$ws.Range("C103").Value = 'Day/'
$ws.Range("D103").Value = 2000
$ws.Range("E103").Value = 2016
$ws.Range("F103").Value = 1700
$ws.Range("G103").Value = 16
$ws.Range("H103").Value = 66
$ws.Range("I103").Value = 200
$ws.Range("J103").Value = 27200
$ws.Range("K103").Value = 'Engine oil 	        250'
$ws.Range("L103").Value = 1200
$ws.Range("M103").Value = 12800
$ws.Range("N103").Value = 'oo'

# Row 104
$ws.Range("A104").Value = "'01-03-2018"
$ws.Range("A104").Style = "Normal"
$ws.Range("B104").Value = 'v4'
$ws.Range("C104").Value = 'Day/'
$ws.Range("D104").Value = 2000
$ws.Range("E104").Value = 2270
$ws.Range("F104").Value = 1700
$ws.Range("G104").Value = 270
$ws.Range("H104").Value = 66
$ws.Range("I104").Value = 200
$ws.Range("J104").Value = 459000
$ws.Range("K104").Value = 'Nothing'
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = 445800
$ws.Range("N104").Value = 'oo'

# Row 105
$ws.Range("A105").Value = "'21-03-2018"
$ws.Range("A105").Style = "Normal"
$ws.Range("B105").Value = 'v1'
$ws.Range("C105").Value = '/Night'
$ws.Range("D105").Value = 1600
$ws.Range("E105").Value = 1624
$ws.Range("F105").Value = 1700
$ws.Range("G105").Value = 24
$ws.Range("H105").Value = 66
$ws.Range("I105").Value = 220
$ws.Range("J105").Value = 40800
$ws.Range("K105").Value = 'Engine oil 	        250'
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 24780

# Row 106
$ws.Range("A106").Value = "'21-03-2018"
$ws.Range("A106").Style = "Normal"
$ws.Range("B106").Value = 'v1'
$ws.Range("C106").Value = '/Night'
$ws.Range("D106").Value = 1600
$ws.Range("E106").Value = 1624
$ws.Range("F106").Value = 1700
$ws.Range("G106").Value = 24
$ws.Range("H106").Value = 66
$ws.Range("I106").Value = 220
$ws.Range("J106").Value = 40800
$ws.Range("K106").Value = 'Engine oil 	        250'
$ws.Range("L106").Value = 1500
$ws.Range("M106").Value = 24780

# Row 107
$ws.Range("A107").Value = "'19-03-2018"
$ws.Range("A107").Style = "Normal"
$ws.Range("B107").Value = 'v4'
$ws.Range("C107").Value = 'Day/Night'
$ws.Range("D107").Value = 2000
$ws.Range("E107").Value = 2016
$ws.Range("F107").Value = 3000
$ws.Range("G107").Value = 16
$ws.Range("H107").Value = 66
$ws.Range("I107").Value = 200
$ws.Range("J107").Value = 48000
$ws.Range("K107").Value = 'Nothing'
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 34800

# Row 108
$ws.Range("A108").Value = "'19-03-2018"
$ws.Range("A108").Style = "Normal"
$ws.Range("B108").Value = 'v4'
$ws.Range("C108").Value = 'Day/Night'
$ws.Range("D108").Value = 2000
$ws.Range("E108").Value = 2016
$ws.Range("F108").Value = 3000
$ws.Range("G108").Value = 16
$ws.Range("H108").Value = 66
$ws.Range("I108").Value = 200
$ws.Range("J108").Value = 48000
$ws.Range("K108").Value = 'Hydraulic strainer  250'
$ws.Range("L108").Value = 1500
$ws.Range("M108").Value = 33300

# Row 109
$ws.Range("A109").Value = "'16-03-2018"
$ws.Range("A109").Style = "Normal"
$ws.Range("B109").Value = 'v3'
$ws.Range("C109").Value = 'Day/Night'
$ws.Range("D109").Value = 2986
$ws.Range("E109").Value = 3000
$ws.Range("F109").Value = 3000
$ws.Range("G109").Value = 14
$ws.Range("H109").Value = 66
$ws.Range("I109").Value = 200
$ws.Range("J109").Value = 42000
$ws.Range("K109").Value = 'Nothing'
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = 28800

# Row 110
$ws.Range("A110").Value = "'16-03-2018"
$ws.Range("A110").Style = "Normal"
$ws.Range("B110").Value = 'v3'
$ws.Range("C110").Value = 'Day/Night'
$ws.Range("D110").Value = 2986
$ws.Range("E110").Value = 3000
$ws.Range("F110").Value = 3000
$ws.Range("G110").Value = 14
$ws.Range("H110").Value = 66
$ws.Range("I110").Value = 200
$ws.Range("J110").Value = 42000
$ws.Range("K110").Value = 'Hydraulic oil	        1000'
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 26800
